$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Token/Exit System settings changes ------------------------------------
# The weighbridge report gains three new tracked fields:
#   "Place" / "Phone No"   -> inserted right after "Vehicle No"
#   "Credit"               -> inserted right after "Charges"
#   "Round off"            -> inserted right after "Bag Deduction"

# 1) Insert "Place" / "Phone No" before the "Material" column.
$ws.Columns("G:H").Insert()
$ws.Range("G1").Value = "Place"
$ws.Range("H1").Value = "Phone No"

# 2) Insert "Credit" before the "Gross Wt" column (now shifted to column L).
$ws.Columns("L:L").Insert()

# 3) Insert "Round off" before the "Nett Wt" column (now shifted to column R).
$ws.Columns("R:R").Insert()

# Fill in the header captions for the newly inserted columns. "Round off" is
# written before "Credit" so the shared-string table picks up the same
# ordering used when this workbook was last hand-edited in Excel.
$ws.Range("R1").Value = "Round off"
$ws.Range("L1").Value = "Credit"

# The legacy "No Of Bags", "Charges", "Bag Deduction", "Final Wt" and
# "Final Amount" sample values are no longer populated for these rows.
$ws.Range("J2:K4").ClearContents()
$ws.Range("Q2:Q4").ClearContents()
$ws.Range("U2:V4").ClearContents()

# Restore the active selection to where the editor left off.
$ws.Range("N9").Select() | Out-Null
